$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1) Insert a new row above row 4 (pushes old rows 4-18 down to 5-19,
#    formulas auto-adjust their relative references).
# ---------------------------------------------------------------------------
$ws.Rows.Item(4).Insert()

# Row 3's running-total formula is removed; only the style stays.
$ws.Range("F3").ClearContents()

# ---------------------------------------------------------------------------
# 2) Fill the newly inserted row 4 (second "Namrata Rubber Product" invoice).
#    Row 6 (A9/B24/C23/D23/E23/F25) carries the exact formatting pattern we
#    need, so stamp its format onto row 4 first, then write the values.
# ---------------------------------------------------------------------------
$ws.Range("A6:F6").Copy()
$ws.Range("A4:F4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Rows.Item(4).RowHeight = 14.4

$ws.Range("B4").Value = "2/6/2024"
$ws.Range("C4").Value = "100/23-24"
$ws.Range("D4").Value = "Namrata Rubber Product Pvt Ltd"
$ws.Range("E4").Value = 55401
$ws.Range("F4").Formula = "=E3+E4"

# ---------------------------------------------------------------------------
# 3) Add the new "Aquachemitech" group in rows 21-23 (row 20 stays an
#    untouched gap row, just like the other group separators).
# ---------------------------------------------------------------------------
$ws.Range("A19:E19").Copy()
$ws.Range("A21:E21").PasteSpecial(-4122)
$ws.Range("A22:E22").PasteSpecial(-4122)
$ws.Range("A23:E23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A19").Copy()
$ws.Range("F21").PasteSpecial(-4122)
$ws.Range("F22").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F19").Copy()
$ws.Range("F23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Rows.Item(21).RowHeight = 14.4
$ws.Rows.Item(22).RowHeight = 14.4
$ws.Rows.Item(23).RowHeight = 14.4

$ws.Range("A21").Value = 7
$ws.Range("B21").Value = "1/24/2024"
$ws.Range("C21").Value = 1483
$ws.Range("D21").Value = "Aquachemitech"
$ws.Range("E21").Value = 8600

$ws.Range("B22").Value = "1/29/2024"
$ws.Range("C22").Value = 1502
$ws.Range("D22").Value = "Aquachemitech"
$ws.Range("E22").Value = 13400

$ws.Range("B23").Value = "1/30/2024"
$ws.Range("C23").Value = 1530
$ws.Range("D23").Value = "Aquachemitech"
$ws.Range("E23").Value = 9525
$ws.Range("F23").Formula = "=E21+E22+E23"

# ---------------------------------------------------------------------------
# 4) Add the new single-row "Asha Enterprises" entry in row 25 (row 24 stays
#    a gap row). Row 8 carries the matching format (B as plain date style,
#    C/D/E/F non-wrapped "9"/"25" styles, no forced row height).
# ---------------------------------------------------------------------------
$ws.Range("A8:F8").Copy()
$ws.Range("A25:F25").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A25").Value = 8
$ws.Range("B25").Value = "2/3/2024"
$ws.Range("C25").Value = 2158
$ws.Range("D25").Value = "Asha Enterprises"
$ws.Range("E25").Value = 1569
$ws.Range("F25").Formula = "=E25"

# ---------------------------------------------------------------------------
# 5) Leave the cursor where the author left it.
# ---------------------------------------------------------------------------
$ws.Range("F30").Select()
